# Add new columns I (I0) and J (IF) with data for rows 2-37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same formatting used by the other header cells (e.g. H1) to the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for rows 2..37: columns I and J
$values = @(
    @(8,8),
    @(7,8),
    @(5,6),
    @(6,6),
    @(9,9),
    @(6,6),
    @(8,8),
    @(8,8),
    @(10,10),
    @(7,7),
    @(12,12),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(6,7),
    @(6,6),
    @(6,7),
    @(6,7),
    @(7,7),
    @(7,7),
    @(4,5),
    @(9,9),
    @(7,8),
    @(4,4),
    @(5,5),
    @(8,8),
    @(9,9),
    @(6,6)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
